# Applies the "demand_parameters.xlsx" edit:
#  - Remove the Mombasa/NH3 row (row 3)
#  - Replace the remaining demand center (Nairobi -> Luderitz) with new
#    coordinates and an updated annual demand figure
#  - Refresh number formatting (comma-separated thousands) and column
#    widths to match the new layout
#  - Move the active selection to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the second demand center (was "Mombasa", row 3) ---------------
$ws.Rows("3:3").Delete() | Out-Null

# --- Update the remaining demand center's data -----------------------------
$ws.Range("A2").Value = "Lüderitz"
$ws.Range("B2").Value = -26.642877645011101
$ws.Range("C2").Value = 15.1439290700957
$ws.Range("D2").Value = 54000000

# --- Number formatting ------------------------------------------------------
# Header for the annual-demand column gets a thousands-separator style.
$ws.Range("D1").Style = "Comma"

# Lat/Lon values use two decimal places.
$ws.Range("B2:C2").NumberFormat = "0.00"

# Annual demand value uses a thousands separator with no decimals.
$ws.Range("D2").Style = "Comma"
$ws.Range("D2").NumberFormat = '_(* #,##0_);_(* \(#,##0\);_(* "-"??_);_(@_)'

# --- Column widths (characters) to fit the new content ----------------------
$ws.Columns("A:A").ColumnWidth = 40.5
$ws.Columns("B:B").ColumnWidth = 15.5
$ws.Columns("C:C").ColumnWidth = 19.333333333333332
$ws.Columns("D:D").ColumnWidth = 21.5
$ws.Columns("E:E").ColumnWidth = 28.833333333333332

# --- Selection moves to the data row now that Mombasa's row is gone --------
$ws.Range("A2").Select() | Out-Null
